# ffmpeg_automate_filename.xlsx edit:
#  - Remove the old helper column S (which just mirrored column C).
#  - Add a new column B holding the path to the ffmpeg executable, used by
#    the command-building formulas in column C.
#  - Update the formulas in column C to reference column B instead of the
#    hard-coded "ffmpeg" command name, and add "-y" (overwrite) flags to
#    both ffmpeg invocations.
#  - Leave the active selection on C17, matching the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused column S (S1:S17 duplicated C1:C17).
$ws.Columns("S:S").Delete()

# New column B: ffmpeg executable path, same value for every data row.
$ws.Range("B1:B17").Value = "/mnt/c/ffmpeg/bin/ffmpeg.exe"

# Row 1 gets its own (non-shared) formula, matching the original layout.
$ws.Range("C1").Formula = '=_xlfn.CONCAT(B1, " -i ", A1,".mp4 -c copy -f h264 -y ", A1, ".h264 && ", B1, " -r 30 -i ", A1, ".h264 -c copy -y ../30fps_input/", A1, "_30fps_input.mp4")'

# Rows 2-17 share one formula (Excel will record it as a shared formula
# rooted at C2, same as before the edit).
$ws.Range("C2:C17").Formula = '=_xlfn.CONCAT(B2, " -i ", A2,".mp4 -c copy -f h264 -y ", A2, ".h264 && ", B2, " -r 30 -i ", A2, ".h264 -c copy -y ../30fps_input/", A2, "_30fps_input.mp4")'

# Match the saved selection/active cell from the target workbook.
$ws.Range("C17").Select()
